$wb = $excel.ActiveWorkbook

# --- Sheet 1 "PI hours": add a new "cfop" column (G) -----------------------
$ws1 = $wb.Worksheets.Item("PI hours")

$ws1.Range("G1").Value = "cfop"
$ws1.Range("G2").Value = "['cfop_MITRA']"
$ws1.Range("G3").Value = "['cfop_NH']"

# Copy the header formatting (bold / border / centered) from an existing
# header cell onto the new header cell.
$ws1.Range("F1").Copy() | Out-Null
$ws1.Range("G1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$excel.CutCopyMode = $false

# --- New sheet "cfop hours" -------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$cfop = $wb.Worksheets.Add($null, $lastSheet)
$cfop.Name = "cfop hours"

$cfop.Range("B1").Value = "cfop"
$cfop.Range("C1").Value = "hours"
$cfop.Range("D1").Value = "percentage"

$cfop.Range("A2").Value = 0
$cfop.Range("B2").Value = "cfop_MITRA"
$cfop.Range("C2").Value = 12
$cfop.Range("D2").Value = 75

$cfop.Range("A3").Value = 1
$cfop.Range("B3").Value = "cfop_NH"
$cfop.Range("C3").Value = 4
$cfop.Range("D3").Value = 25

# Match the look of the other sheets: bold/bordered/centered header row and
# a styled index column (A).
$ws1.Range("B1:D1").Copy() | Out-Null
$cfop.Range("B1:D1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$ws1.Range("A2:A3").Copy() | Out-Null
$cfop.Range("A2:A3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Restore the originally-active sheet/tab selection.
$ws1.Activate()
